$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cobertura de la Prueba")

# Update the "Hs Reales Utilizadas" value for sprint 1 (D3) from 22.5 to 69
$ws.Range("D3").Value = 69

# Move active cell selection to D4 to match post-edit state
$ws.Range("D4").Select()
